$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.306.35'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.082.13'
$ws.Range('E3').Value = '  +3.24%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('D5').Value = '328.18'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').Value = '0.5194'
$ws.Range('E7').Value = '  +1.53%  '
$ws.Range('D8').Value = '0.4312'
$ws.Range('E8').Value = '  +3.53%  '
$ws.Range('D9').Value = '0.08811'
$ws.Range('E9').Value = '  +0.64%  '
$ws.Range('D10').Value = '46.08'
$ws.Range('E10').Value = '  +5.96%  '
$ws.Range('D11').Value = '1.160'
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('D12').Value = '24.50'
$ws.Range('E12').Value = '  -0.61%  '
$ws.Range('D13').Value = '2.080.09'
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('D14').Value = '6.688'
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').Value = '7.682'
$ws.Range('E15').Value = '  +2.63%  '
$ws.Range('D16').Value = '95.24'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '0.00001119'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').Value = '0.06621'
$ws.Range('E19').Value = '  +1.39%  '
$ws.Range('D20').Value = '18.82'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').Value = '0.9989'
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').Value = '6.334'
$ws.Range('E22').Value = '  +1.79%  '
$ws.Range('D23').Value = '30.347.79'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '12.31'
$ws.Range('E24').Value = '  +3.64%  '
$ws.Range('D25').Value = '2.290'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').Value = '2.323.61'
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('D27').Value = '22.31'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').Value = '2.602'
$ws.Range('E28').Value = '  +6.75%  '
$ws.Range('D29').Value = '162.02'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '131.01'
$ws.Range('D31').Value = '1.185'
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('E32').Value = '  +1.39%  '
$ws.Range('D33').Value = '1.631'
$ws.Range('E33').Value = '  +19.27%  '
$ws.Range('D34').Value = '6.202'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('D35').Value = '3.823'
$ws.Range('E35').Value = '  -0.29%  '
$ws.Range('D36').Value = '0.02582'
$ws.Range('E36').Value = '  +2.22%  '
$ws.Range('D37').Value = '9.843'
$ws.Range('E37').Value = '  +8.21%  '
$ws.Range('D38').Value = '12.72'
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('D39').Value = '0.06676'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').Value = '5.446'
$ws.Range('E40').Value = '  -0.83%  '
$ws.Range('D41').Value = '0.2248'
$ws.Range('E41').Value = '  +2.22%  '
$ws.Range('D42').Value = '0.6824'
$ws.Range('E42').Value = '  +2.12%  '
$ws.Range('D43').Value = '1.246'
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.6351'
$ws.Range('E45').Value = '  +2.53%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '13.93'
$ws.Range('E46').Value = '  +1.17%  '
$ws.Range('D47').Value = '2.201'
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('D48').Value = '3.609'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').Value = '1.245'
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('D50').Value = '1.189'
$ws.Range('E50').Value = '  +7.09%  '
$ws.Range('E51').Value = '  +0.43%  '
